$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (number of people interested) column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1491
$ws1.Range("F3").Value = 3146
$ws1.Range("F5").Value = 886
$ws1.Range("F6").Value = 297

# Sheet "全部类型" (all types) - same updates, mapped to the merged rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1491
$ws4.Range("F3").Value = 3146
$ws4.Range("F5").Value = 886
$ws4.Range("F7").Value = 297
